# Apply stimulus updates:
#  1. Column L (correct_ans): recode single-letter answer codes to full words
#     b -> center, y -> left, r -> right
#  2. Columns A-D (promptFile/correctFile/dist_01File/dist_02File): rename the
#     "face" image category to "book" (face//face_NN.jpg -> book//book_NN.jpg)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$answerMap = @{ "b" = "center"; "y" = "left"; "r" = "right" }

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$lastCol = $usedRange.Columns.Count

# --- 1. Recode correct_ans column (column L / 12) ---
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 12)
    $old = $cell.Value2
    if ($answerMap.ContainsKey($old)) {
        $cell.Value = $answerMap[$old]
    }
}

# --- 2. Rename face//face_NN.jpg -> book//book_NN.jpg in columns A-D ---
for ($row = 1; $row -le $lastRow; $row++) {
    for ($col = 1; $col -le 4; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $val = $cell.Value2
        if ($val -like "face//face_*") {
            $cell.Value = ($val -replace "face//face_", "book//book_")
        }
    }
}
